$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the numeric values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0

# Set the text value (goes into shared strings)
$ws.Range("B2").Value = "disconnected_elements"

# Apply formatting (bold font, thin border all around, centered horizontal, top vertical)
# to B1 first (this builds the new font/border/style entries once).
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop
$r1.Borders.LineStyle = 1        # xlContinuous
$r1.Borders.Weight = 2           # xlThin

# Propagate the exact same style to A2 via copy/paste-special (format painter),
# instead of re-applying each property individually, to avoid generating
# spurious intermediate/unused style entries.
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)  # xlPasteFormats
